$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 8 new blank rows (93-100) below the existing data, inheriting
#     the same style (s=12) that rows 13-92 already use ---
$ws.Rows("93:100").Insert()

# --- Fill in the repeating "nha thau" (contractor) sample-data block for
#     rows 6-12, extending the existing pattern already present in rows 4-5 ---
$ws.Range("A6:H12").NumberFormat = "@"

# Row 6 - Nha thau 3
$ws.Cells.Item(6,1).Value  = "NT003"
$ws.Cells.Item(6,2).Value  = "Nhà thầu 3"
$ws.Cells.Item(6,3).Value  = "01-01-2025"
$ws.Cells.Item(6,4).Value  = "7.000.000"
$ws.Cells.Item(6,5).Value  = "01-01-2027"
$ws.Cells.Item(6,6).Value  = "3.000.000"
$ws.Cells.Item(6,7).Value  = "5.000.000"
$ws.Cells.Item(6,8).Value  = "700.000"
$ws.Cells.Item(6,10).Value = 3

# Row 7 - Nha thau 4
$ws.Cells.Item(7,1).Value  = "NT004"
$ws.Cells.Item(7,2).Value  = "Nhà thầu 4"
$ws.Cells.Item(7,3).Value  = "01-01-2025"
$ws.Cells.Item(7,4).Value  = "5.000.000"
$ws.Cells.Item(7,5).Value  = "01-01-2025"
$ws.Cells.Item(7,6).Value  = "1.000.000"
$ws.Cells.Item(7,7).Value  = "3.000.000"
$ws.Cells.Item(7,8).Value  = "500.000"
$ws.Cells.Item(7,9).Value  = 10
$ws.Cells.Item(7,10).Value = 4

# Row 8 - Nha thau 5
$ws.Cells.Item(8,1).Value  = "NT005"
$ws.Cells.Item(8,2).Value  = "Nhà thầu 5"
$ws.Cells.Item(8,3).Value  = "01-01-2025"
$ws.Cells.Item(8,4).Value  = "6.000.000"
$ws.Cells.Item(8,5).Value  = "01-01-2026"
$ws.Cells.Item(8,6).Value  = "2.000.000"
$ws.Cells.Item(8,7).Value  = "4.000.000"
$ws.Cells.Item(8,8).Value  = "600.000"
$ws.Cells.Item(8,10).Value = 5

# Row 9 - Nha thau 6
$ws.Cells.Item(9,1).Value  = "NT006"
$ws.Cells.Item(9,2).Value  = "Nhà thầu 6"
$ws.Cells.Item(9,3).Value  = "01-01-2025"
$ws.Cells.Item(9,4).Value  = "7.000.000"
$ws.Cells.Item(9,5).Value  = "01-01-2027"
$ws.Cells.Item(9,6).Value  = "3.000.000"
$ws.Cells.Item(9,7).Value  = "5.000.000"
$ws.Cells.Item(9,8).Value  = "700.000"
$ws.Cells.Item(9,10).Value = 6

# Row 10 - Nha thau 7
$ws.Cells.Item(10,1).Value  = "NT007"
$ws.Cells.Item(10,2).Value  = "Nhà thầu 7"
$ws.Cells.Item(10,3).Value  = "01-01-2025"
$ws.Cells.Item(10,4).Value  = "5.000.000"
$ws.Cells.Item(10,5).Value  = "01-01-2025"
$ws.Cells.Item(10,6).Value  = "1.000.000"
$ws.Cells.Item(10,7).Value  = "3.000.000"
$ws.Cells.Item(10,8).Value  = "500.000"
$ws.Cells.Item(10,9).Value  = 10
$ws.Cells.Item(10,10).Value = 7

# Row 11 - Nha thau 8
$ws.Cells.Item(11,1).Value  = "NT008"
$ws.Cells.Item(11,2).Value  = "Nhà thầu 8"
$ws.Cells.Item(11,3).Value  = "01-01-2025"
$ws.Cells.Item(11,4).Value  = "6.000.000"
$ws.Cells.Item(11,5).Value  = "01-01-2026"
$ws.Cells.Item(11,6).Value  = "2.000.000"
$ws.Cells.Item(11,7).Value  = "4.000.000"
$ws.Cells.Item(11,8).Value  = "600.000"
$ws.Cells.Item(11,10).Value = 8

# Row 12 - Nha thau 9
$ws.Cells.Item(12,1).Value  = "NT009"
$ws.Cells.Item(12,2).Value  = "Nhà thầu 9"
$ws.Cells.Item(12,3).Value  = "01-01-2025"
$ws.Cells.Item(12,4).Value  = "7.000.000"
$ws.Cells.Item(12,5).Value  = "01-01-2027"
$ws.Cells.Item(12,6).Value  = "3.000.000"
$ws.Cells.Item(12,7).Value  = "5.000.000"
$ws.Cells.Item(12,8).Value  = "700.000"
$ws.Cells.Item(12,10).Value = 9

# --- Extend the data validation on columns I:J to cover the new rows ---
$ws.Range("I4:J100").Validation.Delete()
$ws.Range("I4:J100").Validation.Add(1, 1, 7, "0")
$ws.Range("I4:J100").Validation.IgnoreBlank = $true
$ws.Range("I4:J100").Validation.ShowError = $true
$ws.Range("I4:J100").Validation.ErrorTitle = "Lỗi nhập liệu"
$ws.Range("I4:J100").Validation.ErrorMessage = "Chỉ được nhập số nguyên lớn hơn hoặc bằng 0!"
$ws.Range("I4:J100").Validation.ShowInput = $false

# --- Update the selection to match the edited region ---
$ws.Range("A6:XFD12").Select()
